# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F
$updates = @{
    2  = 35
    3  = 6379
    4  = 185
    7  = 1915
    8  = 1450
    9  = 299
    10 = 975
    11 = 284
    12 = 5593
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
